$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $cellRef, $value) {
    $range = $ws.Range($cellRef)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2
Set-TextValue $ws 'D2' '64.572.80'
$ws.Range('E2').Value = '  +0.71%  '

# Row 3
Set-TextValue $ws 'D3' '3.142.94'
$ws.Range('E3').Value = '  +3.01%  '

# Row 4
Set-TextValue $ws 'D4' '1.00'
$ws.Range('E4').Value = '  +0.16%  '

# Row 5
Set-TextValue $ws 'D5' '564.57'
$ws.Range('E5').Value = '  +2.77%  '

# Row 6
Set-TextValue $ws 'D6' '145.47'
$ws.Range('E6').Value = '  +5.11%  '

# Row 7
Set-TextValue $ws 'D7' '1.00'
$ws.Range('E7').Value = '  +0.20%  '

# Row 8
Set-TextValue $ws 'D8' '3.134.82'
$ws.Range('E8').Value = '  +2.71%  '

# Row 9
Set-TextValue $ws 'D9' '0.496'
$ws.Range('E9').Value = '  +2.11%  '

# Row 10
Set-TextValue $ws 'D10' '6.73'
$ws.Range('E10').Value = '  +4.40%  '

# Row 11
$ws.Range('E11').Value = '  +1.72%  '

# Row 12
Set-TextValue $ws 'D12' '0.468'
$ws.Range('E12').Value = '  +2.53%  '

# Row 13
Set-TextValue $ws 'D13' '36.97'
$ws.Range('E13').Value = '  +3.84%  '

# Row 14
Set-TextValue $ws 'D14' '0.0000223'
$ws.Range('E14').Value = '  +2.07%  '

# Row 15
Set-TextValue $ws 'D15' '3.649.56'
$ws.Range('E15').Value = '  +3.26%  '

# Row 16
Set-TextValue $ws 'D16' '64.631.79'
$ws.Range('E16').Value = '  +0.85%  '

# Row 17
$ws.Range('E17').Value = '  +1.36%  '

# Row 18
Set-TextValue $ws 'D18' '3.140.98'
$ws.Range('E18').Value = '  +3.15%  '

# Row 19
Set-TextValue $ws 'D19' '514.87'
$ws.Range('E19').Value = '  +6.22%  '

# Row 20
Set-TextValue $ws 'D20' '6.86'
$ws.Range('E20').Value = '  +4.36%  '

# Row 21
Set-TextValue $ws 'D21' '14.06'
$ws.Range('E21').Value = '  +3.38%  '

# Row 22
$ws.Range('E22').Value = '  +5.27%  '

# Row 23
Set-TextValue $ws 'D23' '7.47'
$ws.Range('E23').Value = '  +4.45%  '

# Row 24
Set-TextValue $ws 'D24' '12.91'
$ws.Range('E24').Value = '  +4.33%  '

# Row 25
Set-TextValue $ws 'D25' '79.02'
$ws.Range('E25').Value = '  +1.36%  '

# Row 26
Set-TextValue $ws 'D26' '0.996'
$ws.Range('E26').Value = '  -0.38%  '

# Row 27
Set-TextValue $ws 'D27' '8.97'
$ws.Range('E27').Value = '  +16.60%  '

# Row 28
$ws.Range('E28').Value = '  +5.37%  '

# Row 29
$ws.Range('E29').Value = '  +4.52%  '

# Row 30
Set-TextValue $ws 'D30' '1.00'
$ws.Range('E30').Value = '  -0.14%  '

# Row 31
Set-TextValue $ws 'D31' '26.64'
$ws.Range('E31').Value = '  +3.47%  '

# Row 32
$ws.Range('E32').Value = '  +0.10%  '

# Row 33
$ws.Range('E33').Value = '  +2.53%  '

# Row 34
Set-TextValue $ws 'D34' '550.83'
$ws.Range('E34').Value = '  -5.64%  '

# Row 35
Set-TextValue $ws 'D35' '5.41'
$ws.Range('E35').Value = '  +0.27%  '

# Row 36
Set-TextValue $ws 'D36' '6.10'
$ws.Range('E36').Value = '  +3.88%  '

# Row 37
Set-TextValue $ws 'D37' '54.00'
$ws.Range('E37').Value = '  +4.36%  '

# Row 38
$ws.Range('E38').Value = '  +7.28%  '

# Row 39
Set-TextValue $ws 'D39' '0.0827'
$ws.Range('E39').Value = '  +4.62%  '

# Row 40
Set-TextValue $ws 'D40' '3.157.10'
$ws.Range('E40').Value = '  +8.22%  '

# Row 41
$ws.Range('E41').Value = '  +3.79%  '

# Row 42
$ws.Range('B42').Value = 'Cosmos'
$ws.Range('C42').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextValue $ws 'D42' '8.30'
$ws.Range('E42').Value = '  +1.47%  '

# Row 43
$ws.Range('B43').Value = 'dogwifhat'
$ws.Range('C43').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextValue $ws 'D43' '2.73'
$ws.Range('E43').Value = '  -2.06%  '

# Row 44
Set-TextValue $ws 'D44' '0.268'
$ws.Range('E44').Value = '  +10.92%  '

# Row 45
$ws.Range('E45').Value = '  +7.76%  '

# Row 46
Set-TextValue $ws 'D46' '0.999'
$ws.Range('E46').Value = '  +0.02%  '

# Row 47
Set-TextValue $ws 'D47' '25.59'
$ws.Range('E47').Value = '  +3.72%  '

# Row 48
Set-TextValue $ws 'D48' '120.87'
$ws.Range('E48').Value = '  +2.15%  '

# Row 49
Set-TextValue $ws 'D49' '0.0₃0523'
$ws.Range('E49').Value = '  -0.74%  '

# Row 50
$ws.Range('E50').Value = '  +0.53%  '

# Row 51
Set-TextValue $ws 'D51' '2.11'
$ws.Range('E51').Value = '  +4.32%  '
